# Apply the two edits described by the diff:
# 1. Replace the title text "2.2 - Debate I" with "Placeholder - Check Back Later".
# 2. Remove the trailing " :::" runs after "...general edification later."

$d = $word.ActiveDocument

# --- Change 1: Title text -------------------------------------------------
$d.Content.Find.Execute(
    "2.2 - Debate I", $true, $false, $false, $false, $false,
    $true, 1, $false, "Placeholder - Check Back Later", 2
)

# --- Change 2: Drop the trailing " :::" text at the end of the Additional
#     Resources sentence (spans two runs: " " and ":::").
$d.Content.Find.Execute(
    "general edification later. :::", $true, $false, $false, $false, $false,
    $true, 1, $false, "general edification later.", 2
)
